$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 ("I0" / "IF") should use the same style as the
# existing header row (e.g. H1), so copy its formatting across first.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$rowNums = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79)
$iVals   = @(6,9,8,5,7,8,8,7,8,6,8,4,7,7,7,7,7,9,6,7,8,7,7,7,6,8,6,7,7,7,7,8,7,7,7,6,7,6,5,6,7,7,7,8,6,6,6,7,6,6,7,6,4,8,6,9,6,6,7,7,6,7,7,7,8,6,9,7,6,7,7,6,7,7,8,5,5,7)
$jVals   = @(6,9,8,6,7,8,8,8,8,7,8,5,7,7,7,7,7,9,7,8,8,7,7,8,6,8,6,7,7,7,7,9,8,8,7,7,7,7,6,6,7,7,7,8,7,6,6,7,7,7,8,6,5,8,7,9,6,7,7,8,6,7,8,7,8,6,9,8,6,7,7,6,7,7,8,5,6,7)

for ($k = 0; $k -lt $rowNums.Length; $k++) {
    $r = $rowNums[$k]
    $ws.Cells.Item($r, 9).Value = $iVals[$k]
    $ws.Cells.Item($r, 10).Value = $jVals[$k]
}
